$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 11) for year 2021, continuing the 2012-2020 series
# already present in rows 2-10.
$row = 11

$values = @{
    1  = "2021年"
    2  = 368.11
    3  = 110.19
    4  = 18.19
    6  = 197.3
    7  = 714.23
    8  = 26.41
    9  = 1392.17
    10 = 41.89
    11 = 8763.860000000001
    12 = 65.06999999999999
    13 = 2.61
    14 = 0.43
    15 = 87.72
    16 = 44.64
    17 = 0.85
    18 = 12.9
    19 = 242.53
    20 = 15.27
    21 = 1065
    23 = 16.79
    24 = 104.37
    25 = 4.85
    26 = 544.98
    27 = 119.83
    28 = 33.43
    29 = 5.71
    30 = 88.97
    31 = 149.52
    32 = 849.67
    33 = 473.39
    34 = 174.57
    35 = 424.51
    36 = 3.76
    37 = 161.28
    38 = 49
    39 = 174.49
    40 = 5.92
    41 = 928.12
    42 = 43.18
    43 = 1.9
}

foreach ($col in $values.Keys) {
    $ws.Cells.Item($row, $col).Value = $values[$col]
}

# Columns E (5) and V (22) are present-but-empty (empty text) in every
# prior row (2012-2020). A plain "" assignment collapses to a truly blank
# cell in this engine, so seed them as text first ...
$ws.Cells.Item($row, 5).Value = "'"
$ws.Cells.Item($row, 22).Value = "'"

# ... then copy the formats down from the row above (this also clears the
# stray quote-prefix style the seeding step leaves behind) so the new
# cells match their column's existing empty-text formatting exactly, and
# give the new year label (A11) the same formatting as the other year
# labels in column A (bold, bordered, centered).
$ws.Range("A10:AQ10").Copy() | Out-Null
$ws.Range("A11:AQ11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
